# Append the 2025-04-24 Kaspa buy as a new row (14) below the existing data,
# matching the layout of the other manually-appended rows (10, 12, 13):
# column A holds the date as plain text (not an Excel date serial), and
# columns B:D hold plain numbers, none of the four cells carrying an
# explicit style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(14, 1)

# Force text formatting first so Excel doesn't auto-convert the
# "MM/DD/YYYY"-looking string into a date serial number when the value is
# assigned below.
$dateCell.NumberFormat = "@"
$dateCell.Value = "04/24/2025"
# Reset back to the default "Normal" style so the cell ends up with no
# explicit style index, same as the other inline-string date cells above it.
$dateCell.Style = "Normal"

$ws.Cells.Item(14, 2).Value = 509.2249999999985
$ws.Cells.Item(14, 3).Value = 0.09818842358485963
$ws.Cells.Item(14, 4).Value = 50
